# Update the "具体时间范围" (E column) date-range separator from "-" to " - "
# across every sheet, and bump the (G column price-tier) headcount totals
# (F column) for the specific rows that changed.

$wb = $excel.ActiveWorkbook

for ($s = 1; $s -le $wb.Worksheets.Count; $s++) {
    $ws = $wb.Worksheets.Item($s)
    $ur = $ws.UsedRange
    $rowCount = $ur.Rows.Count

    for ($r = 2; $r -le $rowCount; $r++) {
        $cell = $ws.Cells.Item($r, 5)
        $v = $cell.Value2
        if ($v -ne $null -and $v -like '*-*') {
            $new = $v -replace '-', ' - '
            $cell.Value = $new
        }
    }
}

# Per-sheet F-column (headcount) updates: Sheet index -> Row -> New value
# 1 = 展览, 2 = 演出, 3 = 本地生活, 4 = 全部类型

$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(4, 6).Value = 1129
$ws1.Cells.Item(5, 6).Value = 533
$ws1.Cells.Item(8, 6).Value = 686
$ws1.Cells.Item(13, 6).Value = 163
$ws1.Cells.Item(14, 6).Value = 3138
$ws1.Cells.Item(17, 6).Value = 441
$ws1.Cells.Item(20, 6).Value = 286
$ws1.Cells.Item(22, 6).Value = 115
$ws1.Cells.Item(23, 6).Value = 7
$ws1.Cells.Item(24, 6).Value = 678
$ws1.Cells.Item(27, 6).Value = 971
$ws1.Cells.Item(29, 6).Value = 1603
$ws1.Cells.Item(30, 6).Value = 335

$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(7, 6).Value = 236

$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(2, 6).Value = 395
$ws3.Cells.Item(3, 6).Value = 77

$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(2, 6).Value = 395
$ws4.Cells.Item(5, 6).Value = 1129
$ws4.Cells.Item(8, 6).Value = 77
$ws4.Cells.Item(9, 6).Value = 533
$ws4.Cells.Item(12, 6).Value = 686
$ws4.Cells.Item(18, 6).Value = 163
$ws4.Cells.Item(19, 6).Value = 3138
$ws4.Cells.Item(23, 6).Value = 441
$ws4.Cells.Item(26, 6).Value = 286
$ws4.Cells.Item(29, 6).Value = 115
$ws4.Cells.Item(30, 6).Value = 7
$ws4.Cells.Item(31, 6).Value = 236
$ws4.Cells.Item(34, 6).Value = 678
$ws4.Cells.Item(40, 6).Value = 971
$ws4.Cells.Item(42, 6).Value = 1603
$ws4.Cells.Item(43, 6).Value = 336
